$wb = $excel.ActiveWorkbook

# --- "Draw" sheet: replace row 2 data and remove row 3 ---
$wsDraw = $wb.Worksheets.Item("Draw")

$wsDraw.Range("A2").Value = "30-12-2024 14:00"
$wsDraw.Range("B2").Value = "WORLD"
$wsDraw.Range("C2").Value = "FRIENDLIES CLUBS"
$wsDraw.Range("D2").Value = "Botafogo PB - Serra Branca"
$wsDraw.Range("E2").Value = 60
$wsDraw.Range("F2").Value = 3.35

$wsDraw.Rows.Item(3).Delete()

# --- "Over_Under" sheet: remove old rows 2 & 3, keep former row 4 as new row 2 ---
$wsOverUnder = $wb.Worksheets.Item("Over_Under")

$wsOverUnder.Rows.Item(2).Delete()
$wsOverUnder.Rows.Item(2).Delete()

$wb.Save()
